$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'67.868.73"
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -1.41%  '
$cell = $ws.Range("D3")
$cell.Value = "'3.842.25"
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").Value = '  +0.29%  '
$cell = $ws.Range("D5")
$cell.Value = "'596.12"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.57%  '
$cell = $ws.Range("D6")
$cell.Value = "'166.13"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '
$cell = $ws.Range("D7")
$cell.Value = "'3.842.06"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -2.26%  '
$ws.Range("E8").Value = '  +0.27%  '
$cell = $ws.Range("D9")
$cell.Value = "'0.525"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.36%  '
$cell = $ws.Range("D10")
$cell.Value = "'0.164"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -1.53%  '
$cell = $ws.Range("D11")
$cell.Value = "'6.29"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -2.15%  '
$cell = $ws.Range("D12")
$cell.Value = "'0.457"
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.15%  '
$cell = $ws.Range("D13")
$cell.Value = "'0.0000247"
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.25%  '
$cell = $ws.Range("D14")
$cell.Value = "'36.87"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.19%  '
$cell = $ws.Range("D15")
$cell.Value = "'4.486.12"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -2.10%  '
$cell = $ws.Range("D16")
$cell.Value = "'3.852.28"
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.82%  '
$cell = $ws.Range("D17")
$cell.Value = "'67.931.33"
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.48%  '
$cell = $ws.Range("D18")
$cell.Value = "'7.45"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '
$cell = $ws.Range("D19")
$cell.Value = "'18.04"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +5.53%  '
$cell = $ws.Range("D20")
$cell.Value = "'0.111"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.67%  '
$cell = $ws.Range("D21")
$cell.Value = "'10.73"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -3.90%  '
$cell = $ws.Range("D22")
$cell.Value = "'468.11"
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -4.25%  '
$cell = $ws.Range("D23")
$cell.Value = "'0.729"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +0.74%  '
$cell = $ws.Range("D24")
$cell.Value = "'0.0000160"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -3.62%  '
$cell = $ws.Range("D25")
$cell.Value = "'83.96"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.63%  '
$cell = $ws.Range("D26")
$cell.Value = "'2.20"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -2.83%  '
$cell = $ws.Range("D27")
$cell.Value = "'12.16"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range("D28")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D29")
$cell.Value = "'9.96"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.00%  '
$cell = $ws.Range("D30")
$cell.Value = "'2.91"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.49%  '
$cell = $ws.Range("D31")
$cell.Value = "'3.992.28"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.10%  '
$cell = $ws.Range("D32")
$cell.Value = "'7.70"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -2.52%  '
$cell = $ws.Range("D33")
$cell.Value = "'2.29"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -4.57%  '
$cell = $ws.Range("D34")
$cell.Value = "'30.89"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -4.62%  '
$cell = $ws.Range("D35")
$cell.Value = "'3.815.49"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.58%  '
$cell = $ws.Range("D36")
$cell.Value = "'0.104"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.08%  '
$cell = $ws.Range("D37")
$cell.Value = "'0.139"
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("E38").Value = '  -3.06%  '
$cell = $ws.Range("D39")
$cell.Value = "'5.89"
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$cell = $ws.Range("D40")
$cell.Value = "'3.25"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +7.93%  '
$ws.Range("E41").Value = '  +0.19%  '
$cell = $ws.Range("D42")
$cell.Value = "'0.312"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -3.23%  '
$cell = $ws.Range("D43")
$cell.Value = "'1.98"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.99%  '
$cell = $ws.Range("D44")
$cell.Value = "'423.66"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -4.16%  '
$ws.Range("E45").Value = '  +0.00%  '
$cell = $ws.Range("D46")
$cell.Value = "'47.21"
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -2.64%  '
$cell = $ws.Range("D47")
$cell.Value = "'8.55"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.20%  '
$cell = $ws.Range("D48")
$cell.Value = "'142.98"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.85%  '
$cell = $ws.Range("D49")
$cell.Value = "'0.000267"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +13.17%  '
$cell = $ws.Range("D50")
$cell.Value = "'0.0355"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.84%  '
$cell = $ws.Range("D51")
$cell.Value = "'39.11"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.43%  '